# Update issue 226: SALVIA:   Update performance document.
$wb = $excel.ActiveWorkbook

$wsDebug   = $wb.Worksheets.Item(1)   # "20130105-PartOfSponza-Debug"
$wsRelease = $wb.Worksheets.Item(2)   # "20130105-PartOfSponza-Release"

# --- Rename sheets: drop the "20130105-" date prefix ---
$wsDebug.Name   = "PartOfSponza-Debug"
$wsRelease.Name = "PartOfSponza-Release"

# --- Debug sheet: header for column C now reads "v1231" (was "v1331") ---
$wsDebug.Range("C1").Value = "v1231"

# --- Release sheet: add the new "v1232" column (D) of step data ---
$wsRelease.Range("D1").Value = "v1232"
$wsRelease.Range("D2").Value = 75
$wsRelease.Range("D3").Value = 76
$wsRelease.Range("D4").Value = 76
$wsRelease.Range("D5").Value = 76
$wsRelease.Range("D6").Value = 76
$wsRelease.Range("D7").Value = 75
$wsRelease.Range("D8").Value = 76
$wsRelease.Range("D9").Value = 75
$wsRelease.Range("D10").Value = 75
$wsRelease.Range("D11").Value = 75

# --- Update selections on each sheet ---
$wsDebug.Range("C2").Select()

# Activating the Release sheet last makes it the active/selected tab,
# matching the new tabSelected/activeTab state.
$wsRelease.Activate()
$wsRelease.Range("G3").Select()
